# Atjaunots labojumu excel dok.
# Insert a new tracked-item row before row 23 (continuing the numbered
# list that ends at row 22 with item 21). Everything below row 22 shifts
# down by one row; Excel's row Insert() carries the formatting of the
# row above down into the freshly-inserted row, which matches the style
# pattern (s=16 / s=17) used by the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, pushing the legend/footer rows down.
$ws.Rows("23").Insert()

# New row continues the numbering from row 22 (item 21) -> item 22.
$ws.Range("A23").Value = 22

# Match the author's final selection/active cell.
$ws.Range("B23").Select()
